$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 13640.167
$ws.Range("I18").Value = 7685.4287
$ws.Range("K18").Value = 7685.4287
$ws.Range("M18").Value = -7401.4287

$ws.Range("H62").Value = 1134.5
$ws.Range("I62").Value = 1134.5
$ws.Range("K62").Value = 1134.5
$ws.Range("M62").Value = -510.5

$ws.Range("H65").Value = 1134.5
$ws.Range("I65").Value = 1134.5
$ws.Range("K65").Value = 5672.5
$ws.Range("M65").Value = -2552.5

$ws.Range("H113").Value = 21169.062
$ws.Range("I113").Value = 38876
$ws.Range("J113").Value = 3462.125
$ws.Range("K113").Value = 38876
$ws.Range("L113").Value = 3462.125
$ws.Range("M113").Value = -35622
$ws.Range("N113").Value = -9970.125

$ws.Range("H129").Value = 1281.6086
$ws.Range("I129").Value = 377
$ws.Range("K129").Value = 1131
$ws.Range("M129").Value = 3869

$ws.Range("H135").Value = 747.5
$ws.Range("J135").Value = 1092.4286
$ws.Range("L135").Value = 9831.857399999999
$ws.Range("N135").Value = -14901.8574

$ws.Range("H137").Value = 1651.25
$ws.Range("I137").Value = 1001
$ws.Range("J137").Value = 2735
$ws.Range("K137").Value = 3003
$ws.Range("L137").Value = 8205
$ws.Range("M137").Value = -453
$ws.Range("N137").Value = -13305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1553.7333
$ws.Range("I45").Value = 1040
$ws.Range("J45").Value = 1810.6
$ws.Range("K45").Value = 1040
$ws.Range("L45").Value = 1810.6
$ws.Range("M45").Value = -663
$ws.Range("N45").Value = -2564.6

$ws.Range("H61").Value = 5098.6113
$ws.Range("I61").Value = 4797.423
$ws.Range("J61").Value = 5881.7
$ws.Range("K61").Value = 4797.423
$ws.Range("L61").Value = 5881.7
$ws.Range("M61").Value = -4585.423
$ws.Range("N61").Value = -6305.7

$ws.Range("H74").Value = 888.8108
$ws.Range("I74").Value = 574.625
$ws.Range("K74").Value = 574.625
$ws.Range("M74").Value = 299.375

$ws.Range("H77").Value = 888.8108
$ws.Range("I77").Value = 574.625
$ws.Range("K77").Value = 2873.125
$ws.Range("M77").Value = 1494.875

$ws.Range("H136").Value = 5098.6113
$ws.Range("I136").Value = 4797.423
$ws.Range("J136").Value = 5881.7
$ws.Range("K136").Value = 14392.269
$ws.Range("L136").Value = 17645.1
$ws.Range("M136").Value = -11842.269
$ws.Range("N136").Value = -22745.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws.Range("H134").Value = 4218.2563
$ws.Range("I134").Value = 4348.243
$ws.Range("K134").Value = 13044.729
$ws.Range("M134").Value = -10509.729

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3521.8333
$ws.Range("I31").Value = 5128.6665
$ws.Range("J31").Value = 2986.2222
$ws.Range("K31").Value = 5128.6665
$ws.Range("L31").Value = 2986.2222
$ws.Range("M31").Value = -4833.6665
$ws.Range("N31").Value = -3576.2222

$ws.Range("H34").Value = 3521.8333
$ws.Range("I34").Value = 5128.6665
$ws.Range("J34").Value = 2986.2222
$ws.Range("K34").Value = 5128.6665
$ws.Range("L34").Value = 2986.2222
$ws.Range("M34").Value = -4926.6665
$ws.Range("N34").Value = -3390.2222

$ws.Range("H58").Value = 2418566.2
$ws.Range("I58").Value = 4350968.5
$ws.Range("J58").Value = 3063.75
$ws.Range("K58").Value = 4350968.5
$ws.Range("L58").Value = 3063.75
$ws.Range("M58").Value = -4350765.5
$ws.Range("N58").Value = -3469.75

$ws.Range("H122").Value = 2369.8667
$ws.Range("I122").Value = 2180.75
$ws.Range("J122").Value = 2438.6365
$ws.Range("K122").Value = 6542.25
$ws.Range("L122").Value = 7315.9095
$ws.Range("M122").Value = -4092.25
$ws.Range("N122").Value = -12215.9095

$ws.Range("H134").Value = 1666.4375
$ws.Range("I134").Value = 1517.738
$ws.Range("K134").Value = 4553.214
$ws.Range("M134").Value = -2018.214

$ws.Range("H136").Value = 2418566.2
$ws.Range("I136").Value = 4350968.5
$ws.Range("J136").Value = 3063.75
$ws.Range("K136").Value = 13052905.5
$ws.Range("L136").Value = 9191.25
$ws.Range("M136").Value = -13050355.5
$ws.Range("N136").Value = -14291.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 579.1429000000001
$ws.Range("I5").Value = 492.875
$ws.Range("J5").Value = 855.2
$ws.Range("K5").Value = 1478.625
$ws.Range("L5").Value = 2565.6
$ws.Range("M5").Value = -1366.625
$ws.Range("N5").Value = -2789.6

$ws.Range("H17").Value = 13055.889
$ws.Range("I17").Value = 249.5
$ws.Range("J17").Value = 16714.857
$ws.Range("K17").Value = 748.5
$ws.Range("L17").Value = 50144.571
$ws.Range("M17").Value = -579.5
$ws.Range("N17").Value = -50482.571

$ws.Range("H38").Value = 398.83334
$ws.Range("I38").Value = 97.75
$ws.Range("J38").Value = 1001
$ws.Range("K38").Value = 293.25
$ws.Range("L38").Value = 3003
$ws.Range("M38").Value = 53.75
$ws.Range("N38").Value = -3697

$ws.Range("H80").Value = 1957.6
$ws.Range("I80").Value = 1194
$ws.Range("J80").Value = 2466.6667
$ws.Range("K80").Value = 3582
$ws.Range("L80").Value = 7400.000100000001
$ws.Range("M80").Value = -2646
$ws.Range("N80").Value = -9272.000100000001

$ws.Range("H83").Value = 1957.6
$ws.Range("I83").Value = 1194
$ws.Range("J83").Value = 2466.6667
$ws.Range("K83").Value = 10746
$ws.Range("L83").Value = 22200.0003
$ws.Range("M83").Value = -6066
$ws.Range("N83").Value = -31560.0003

$ws.Range("H98").Value = 1616.6666
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 1616.6666
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 4849.9998
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -7845.9998

$ws.Range("H135").Value = 579.1429000000001
$ws.Range("I135").Value = 492.875
$ws.Range("J135").Value = 855.2
$ws.Range("K135").Value = 4435.875
$ws.Range("L135").Value = 7696.8
$ws.Range("M135").Value = -1900.875
$ws.Range("N135").Value = -12766.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1733.3334
$ws.Range("J122").Value = 2100
$ws.Range("L122").Value = 6300
$ws.Range("N122").Value = -11200

$ws.Range("H132").Value = 1041482
$ws.Range("I132").Value = 1375153.6
$ws.Range("J132").Value = 3392.5557
$ws.Range("K132").Value = 4125460.8
$ws.Range("L132").Value = 10177.6671
$ws.Range("M132").Value = -4122930.8
$ws.Range("N132").Value = -15237.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 16271.637
$ws.Range("J40").Value = 12398.5
$ws.Range("L40").Value = 12398.5
$ws.Range("N40").Value = -12670.5

$ws.Range("H46").Value = 1990.4286
$ws.Range("J46").Value = 3046.5
$ws.Range("L46").Value = 3046.5
$ws.Range("N46").Value = -3422.5

$ws.Range("H61").Value = 3197.2964
$ws.Range("I61").Value = 2891.1
$ws.Range("J61").Value = 4072.1428
$ws.Range("K61").Value = 2891.1
$ws.Range("L61").Value = 4072.1428
$ws.Range("M61").Value = -2689.1
$ws.Range("N61").Value = -4476.1428

$ws.Range("H100").Value = 1097
$ws.Range("I100").Value = 1097
$ws.Range("K100").Value = 1097
$ws.Range("M100").Value = -556

$ws.Range("H113").Value = 3197.2964
$ws.Range("I113").Value = 2891.1
$ws.Range("J113").Value = 4072.1428
$ws.Range("K113").Value = 2891.1
$ws.Range("L113").Value = 4072.1428
$ws.Range("M113").Value = -721.0999999999999
$ws.Range("N113").Value = -8412.1428

$ws.Range("H127").Value = 34389.332
$ws.Range("J127").Value = 34389.332
$ws.Range("L127").Value = 34389.332
$ws.Range("N127").Value = -44309.332

$ws.Range("H136").Value = 1924.76
$ws.Range("I136").Value = 1196.2632
$ws.Range("J136").Value = 4231.6665
$ws.Range("K136").Value = 3588.7896
$ws.Range("L136").Value = 12694.9995
$ws.Range("M136").Value = -1038.7896
$ws.Range("N136").Value = -17794.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 673.1852
$ws.Range("I107").Value = 529.82355
$ws.Range("K107").Value = 1589.47065
$ws.Range("M107").Value = 330.5293500000002

$ws.Range("H128").Value = 35000
$ws.Range("J128").Value = 35000
$ws.Range("L128").Value = 35000
$ws.Range("N128").Value = -44960

$ws.Range("H136").Value = 15433341
$ws.Range("I136").Value = 20577030
$ws.Range("K136").Value = 61731090
$ws.Range("M136").Value = -61728540
